# Append the "first cluster" run result (00000863) to Sheet1 and clear the
# stray empty C2 cell left over from the previous run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: C2 was an empty placeholder cell -> drop it entirely.
$ws.Range("C2").Value = ""

# Row 3: new run data. A3 must keep its leading zeros (text), so force a
# text number format before assigning the value, then reset the style back
# to Normal/General so no extra style index is left on the cell.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "00000863"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = 11615001
$ws.Range("C3").Value = 0.4
$ws.Range("D3").Value = 0.0001
$ws.Range("E3").Value = 0.9
$ws.Range("F3").Value = 0.8
$ws.Range("G3").Value = 300
$ws.Range("H3").Value = 300
